# Weekly update: a new week's price record is inserted at the top of the
# data block (row 74), pushing the existing history down by one row
# (old row 100 ends up duplicated at the new row 101) and the sheet's
# used range grows from A1:R100 to A1:R101.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 74:100 down to 75:101, leaving a blank row 74 to fill in.
$ws.Rows.Item(74).Insert()

# Populate the new row 74 with this week's record (same categorical
# attributes as the series, new date + price figures).
$ws.Range("A74").Value = 9
$ws.Range("B74").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C74").Value = "Metropolitana"
$ws.Range("D74").Value = 44755
$ws.Range("E74").Value = 13
$ws.Range("F74").Value = 100112005
$ws.Range("G74").Value = "Puerro"
$ws.Range("H74").Value = "Sin especificar"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 160
$ws.Range("K74").Value = 7000
$ws.Range("L74").Value = 8000
$ws.Range("M74").Value = 7500
$ws.Range("N74").Value = "$/paquete 20 unidades"
$ws.Range("O74").Value = "Provincia de Chacabuco"
$ws.Range("P74").Value = 375
$ws.Range("Q74").Value = 20
$ws.Range("R74").Value = "Hortaliza"
